$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AJ1").Value = "24-jul"

$ws.Range("AJ2").Value = 0
$ws.Range("AJ3").Value = 9.8723034898639028
$ws.Range("AJ4").Value = 11.383038046902159
$ws.Range("AJ5").Value = 28.803963576257583
$ws.Range("AJ6").Value = 0
$ws.Range("AJ7").Value = 18.669100546592478
$ws.Range("AJ8").Value = 7.7688735866172669
$ws.Range("AJ9").Value = 18.975278738927418
$ws.Range("AJ10").Value = 23.011524730821645
$ws.Range("AJ11").Value = 16.781644240659087
$ws.Range("AJ12").Value = 0
$ws.Range("AJ13").Value = 16.819050389046506
$ws.Range("AJ14").Value = 0
$ws.Range("AJ15").Value = 0
$ws.Range("AJ16").Value = 10.584004969564333
$ws.Range("AJ17").Value = 0
$ws.Range("AJ18").Value = 0

$ws.Range("AJ1").NumberFormat = "@"

$ws.Range("AL8").Select()
